$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new Job Posting row (row 7) with Job_Id = JD_006
$ws.Range("A7").Value = "JD_006"
$ws.Range("B7").Value = "Junior Devops Engineer"
$ws.Range("C7").Value = "Testing"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 3
